$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 657.5
$ws.Range("I12").Value = 636.75
$ws.Range("K12").Value = 636.75
$ws.Range("M12").Value = -466.75
$ws.Range("H17").Value = 2002.7222
$ws.Range("J17").Value = 2002.7222
$ws.Range("L17").Value = 6008.1666
$ws.Range("N17").Value = -6344.1666
$ws.Range("H34").Value = 10006.4
$ws.Range("I34").Value = 10006.4
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 10006.4
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -9803.4
$ws.Range("N34").Value = $null
$ws.Range("H36").Value = 10006.4
$ws.Range("I36").Value = 10006.4
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 10006.4
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9291.4
$ws.Range("N36").Value = $null
$ws.Range("H40").Value = 3533.1667
$ws.Range("I40").Value = 2133
$ws.Range("J40").Value = 4933.3335
$ws.Range("K40").Value = 2133
$ws.Range("L40").Value = 4933.3335
$ws.Range("M40").Value = -1958
$ws.Range("N40").Value = -5283.3335
$ws.Range("H42").Value = 327.5
$ws.Range("I42").Value = 55
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 165
$ws.Range("L42").Value = 1800
$ws.Range("M42").Value = 65
$ws.Range("N42").Value = -2260
$ws.Range("H135").Value = 1083.5
$ws.Range("I135").Value = 1125.25
$ws.Range("K135").Value = 10127.25
$ws.Range("M135").Value = -7592.25
$ws.Range("H138").Value = 2133
$ws.Range("J138").Value = 2898.3333
$ws.Range("L138").Value = 8694.999899999999
$ws.Range("N138").Value = -18974.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 377
$ws.Range("I5").Value = 377
$ws.Range("K5").Value = 377
$ws.Range("M5").Value = -265
$ws.Range("H6").Value = 1001499
$ws.Range("I6").Value = 1001499
$ws.Range("K6").Value = 1001499
$ws.Range("M6").Value = -1001326
$ws.Range("H50").Value = 5326.5
$ws.Range("I50").Value = 2102
$ws.Range("J50").Value = 15000
$ws.Range("K50").Value = 2102
$ws.Range("L50").Value = 15000
$ws.Range("M50").Value = -1388
$ws.Range("N50").Value = -16428
$ws.Range("H135").Value = 69498.5
$ws.Range("J135").Value = 69498.5
$ws.Range("L135").Value = 69498.5
$ws.Range("N135").Value = -79638.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 377
$ws.Range("I4").Value = 377
$ws.Range("K4").Value = 377
$ws.Range("M4").Value = -262
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = $null
$ws.Range("H86").Value = 1171.75
$ws.Range("I86").Value = 1210.5714
$ws.Range("K86").Value = 1210.5714
$ws.Range("M86").Value = -87.57140000000004
$ws.Range("H89").Value = 1171.75
$ws.Range("I89").Value = 1210.5714
$ws.Range("K89").Value = 6052.857
$ws.Range("M89").Value = -436.857
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = $null

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 866.1667
$ws.Range("I5").Value = 274.5
$ws.Range("K5").Value = 274.5
$ws.Range("M5").Value = -162.5
$ws.Range("H22").Value = 1695
$ws.Range("I22").Value = 1695
$ws.Range("K22").Value = 1695
$ws.Range("M22").Value = -1345
$ws.Range("H25").Value = 2414.2856
$ws.Range("I25").Value = 2414.2856
$ws.Range("K25").Value = 2414.2856
$ws.Range("M25").Value = -2240.2856
$ws.Range("H41").Value = 12039.444
$ws.Range("I41").Value = 9231.875
$ws.Range("J41").Value = 34500
$ws.Range("K41").Value = 9231.875
$ws.Range("L41").Value = 34500
$ws.Range("M41").Value = -8803.875
$ws.Range("N41").Value = -35356
$ws.Range("H58").Value = 2127.75
$ws.Range("I58").Value = 2006
$ws.Range("J58").Value = 2249.5
$ws.Range("K58").Value = 2006
$ws.Range("L58").Value = 2249.5
$ws.Range("M58").Value = -1803
$ws.Range("N58").Value = -2655.5
$ws.Range("H136").Value = 2127.75
$ws.Range("I136").Value = 2006
$ws.Range("J136").Value = 2249.5
$ws.Range("K136").Value = 6018
$ws.Range("L136").Value = 6748.5
$ws.Range("M136").Value = -3468
$ws.Range("N136").Value = -11848.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 199.875
$ws.Range("I12").Value = 257.8
$ws.Range("J12").Value = 103.333336
$ws.Range("K12").Value = 773.4000000000001
$ws.Range("L12").Value = 310.000008
$ws.Range("M12").Value = -600.4000000000001
$ws.Range("N12").Value = -656.000008
$ws.Range("H18").Value = 230.22223
$ws.Range("I18").Value = 230.22223
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 690.66669
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -521.66669
$ws.Range("N18").Value = $null
$ws.Range("H50").Value = 2529.6
$ws.Range("I50").Value = 235
$ws.Range("K50").Value = 705
$ws.Range("M50").Value = -224
$ws.Range("H51").Value = 897.5
$ws.Range("I51").Value = 863.3333
$ws.Range("K51").Value = 2589.9999
$ws.Range("M51").Value = -2129.9999
$ws.Range("H53").Value = 2529.6
$ws.Range("I53").Value = 235
$ws.Range("K53").Value = 705
$ws.Range("M53").Value = -224
$ws.Range("H117").Value = 1889
$ws.Range("I117").Value = 2176
$ws.Range("J117").Value = 1028
$ws.Range("K117").Value = 6528
$ws.Range("L117").Value = 3084
$ws.Range("M117").Value = -3086
$ws.Range("N117").Value = -9968

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 216.125
$ws.Range("I2").Value = 216.125
$ws.Range("K2").Value = 216.125
$ws.Range("M2").Value = -103.125
$ws.Range("H11").Value = 18246796
$ws.Range("J11").Value = 99943.5
$ws.Range("L11").Value = 99943.5
$ws.Range("N11").Value = -100221.5
$ws.Range("H21").Value = 62251
$ws.Range("J21").Value = 62251
$ws.Range("L21").Value = 62251
$ws.Range("N21").Value = -62597
$ws.Range("H24").Value = 29748.875
$ws.Range("J24").Value = 29748.875
$ws.Range("L24").Value = 29748.875
$ws.Range("N24").Value = -30094.875
$ws.Range("H30").Value = 62251
$ws.Range("J30").Value = 62251
$ws.Range("L30").Value = 62251
$ws.Range("N30").Value = -62461
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 7500
$ws.Range("M83").Value = -2508
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null
$ws.Range("H55").Value = 889.8461
$ws.Range("I55").Value = 216.85715
$ws.Range("K55").Value = 216.85715
$ws.Range("M55").Value = -43.85714999999999
$ws.Range("H59").Value = 32499.5
$ws.Range("J59").Value = 32499.5
$ws.Range("L59").Value = 32499.5
$ws.Range("N59").Value = -33807.5
$ws.Range("I136").Value = 12109.889
$ws.Range("K136").Value = 36329.667
$ws.Range("M136").Value = -33779.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 9313
$ws.Range("I22").Value = 9313
$ws.Range("K22").Value = 9313
$ws.Range("M22").Value = -9020
$ws.Range("H45").Value = 129999
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("H126").Value = 2799.818
$ws.Range("I126").Value = 2766.25
$ws.Range("K126").Value = 8298.75
$ws.Range("M126").Value = -5828.75
$ws.Range("H132").Value = 6202
$ws.Range("I132").Value = 7001.75
$ws.Range("J132").Value = 3003
$ws.Range("K132").Value = 21005.25
$ws.Range("L132").Value = 9009
$ws.Range("M132").Value = -18475.25
$ws.Range("N132").Value = -14069
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
